$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.007739543914795
$ws.Range("B1").Value = 6.196311950683594
$ws.Range("C1").Value = 3.26241135597229
$ws.Range("D1").Value = 1.438793778419495
$ws.Range("E1").Value = 1.010600328445435
